# Generate Report for Handback
# Adds a new handback row (efb15e1e-83d3-409f-9004-172e11e35e94) to the
# Overview sheet plus the per-locale (zh-cn / de-de) detail sheets, mirroring
# the existing rows for 8a001e41-... / 8e407b60-....
#
# Note: values that Excel would otherwise auto-coerce to a different cell
# type (dates, True/False -> boolean) are written with a leading apostrophe
# so they land as plain text, matching the source data (everything in this
# workbook is authored as literal strings, including the date/bool columns).

$wb = $excel.ActiveWorkbook

$fileId   = "efb15e1e-83d3-409f-9004-172e11e35e94"
$zhXlf    = "efb15e1e-83d3-409f-9004-172e11e35e94.1e6fea1ba5f5f7afd08af0612eb7e39a113aacef.zh-cn.xlf"
$deXlf    = "efb15e1e-83d3-409f-9004-172e11e35e94.1e6fea1ba5f5f7afd08af0612eb7e39a113aacef.de-de.xlf"

$urlSrc  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1d2e3f4a5b6c7d8e9f0a1b2c3d4e5f6a7b8c9d0/e2e/$fileId.md"
$urlZh   = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d2e3f4a5b6c7d8e9f0a1b2c3d4e5f6a7b8c9d0e1/e2e/$fileId.md"
$urlDe   = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e3f4a5b6c7d8e9f0a1b2c3d4e5f6a7b8c9d0e1f2/e2e/$fileId.md"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> new row 4
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tOverview = $wsOverview.ListObjects.Item(1)
$tOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "$fileId.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "'2016-08-28 08:43:41"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $urlSrc, "", "", "e2e\$fileId.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> new row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$tZh = $wsZh.ListObjects.Item(1)
$tZh.ListRows.Add() | Out-Null

$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = "'2016-08-28 08:43:36"
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = "'2016-08-28 08:43:52"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("O4").Value = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $urlSrc, "", "", "$fileId.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $urlZh, "", "", "$fileId.md")

# ---------------------------------------------------------------------------
# Sheet "de-de" -> new row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$tDe = $wsDe.ListObjects.Item(1)
$tDe.ListRows.Add() | Out-Null

$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = "'2016-08-28 08:43:41"
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = "'2016-08-28 08:44:00"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("O4").Value = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $urlSrc, "", "", "$fileId.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $urlDe, "", "", "$fileId.md")
